$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.833.89'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '2.484.50'
$ws.Range('E3').Value = '  -1.70%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('E5').Value = '  -1.68%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '104.04'
$ws.Range('E6').Value = '  -5.00%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.518'
$ws.Range('E7').Value = '  -2.93%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.534'
$ws.Range('E9').Value = '  -3.44%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '38.59'
$ws.Range('E10').Value = '  -4.76%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '20.46'
$ws.Range('E11').Value = '  -0.54%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0798'
$ws.Range('E12').Value = '  -3.33%  '
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('E14').Value = '  -3.94%  '
$ws.Range('D15').Value = '2.871.75'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('D16').Value = '2.483.14'
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.824'
$ws.Range('E17').Value = '  -3.68%  '
$ws.Range('D18').Value = '47.754.60'
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('E19').Value = '  +7.76%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.67'
$ws.Range('E20').Value = '  -6.51%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.52'
$ws.Range('E21').Value = '  -1.90%  '
$ws.Range('E22').Value = '  -2.66%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '278.55'
$ws.Range('E23').Value = '  +5.18%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '70.70'
$ws.Range('E24').Value = '  -1.83%  '
$ws.Range('E25').Value = '  -3.35%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '25.66'
$ws.Range('E27').Value = '  -1.78%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.22'
$ws.Range('E28').Value = '  -6.13%  '
$ws.Range('E29').Value = '  -5.76%  '
$ws.Range('E30').Value = '  -5.36%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '34.37'
$ws.Range('E31').Value = '  -4.26%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '49.14'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '18.90'
$ws.Range('E34').Value = '  -4.60%  '
$ws.Range('E35').Value = '  -3.00%  '
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('E37').Value = '  -3.19%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '4.50'
$ws.Range('E38').Value = '  -4.66%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.85'
$ws.Range('E39').Value = '  -5.34%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '122.99'
$ws.Range('E40').Value = '  +1.45%  '
$ws.Range('E41').Value = '  -1.69%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '21.56'
$ws.Range('E44').Value = '  -1.16%  '
$ws.Range('D45').Value = '1.990.79'
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.11'
$ws.Range('E46').Value = '  -2.37%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.87'
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('E48').Value = '  -4.13%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.89'
$ws.Range('E49').Value = '  -2.65%  '
$ws.Range('E50').Value = '  -3.15%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '78.59'
$ws.Range('E51').Value = '  -0.94%  '
